# Generate Report for handback
# Updates the "Correspond Handoff Datetime" (D3) and
# "Correspond Handback DateTime" (G3) timestamps for the
# 889e6598-... entry on both the zh-cn and de-de language sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-18 05:46:17"
$wsZhCn.Range("G3").Value = "2016-01-18 05:47:06"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-18 05:46:30"
$wsDeDe.Range("G3").Value = "2016-01-18 05:47:30"
